$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Visión (row 8) -> fill in the previously-empty description
$ws.Range("C8").Value = "Ser la marca de óptica reconocida y preferida por los colombianos."

# Misión (row 9) -> fill in the previously-empty description
$ws.Range("C9").Value = "Nuestra misión es ofrecer monturas y lentes de calidad a precios razonables, con diseños característicos y elegantes, que usen materiales resistentes y que sean reconocibles."

# Propósito superior (row 6) -> new text (old text moves down to Objetivo retador)
$ws.Range("C6").Value = "Queremos brindar visión a nuestros clientes, porque ver el mundo te ayuda a comprenderlo."

# Objetivo retador (row 7) -> now holds what used to be the Propósito superior's old 2030-goal text
$ws.Range("C7").Value = "Para el 2030 llegar a 100 ópticas abiertas a nivel nacional."

# Valores (principios) table rows 13-17: fill in Nombre (B) and Descripción (C)
$ws.Range("B13").Value = "Responsabilidad ambiental"
$ws.Range("C13").Value = "En la medida de lo posible nuestros productos son amigables con el medio ambiente, ya que los compramos siendo hechos con materiales reciclados"

$ws.Range("B14").Value = "Transparencia"
$ws.Range("C14").Value = "Nuestra prioridad es la salud del cliente. Todos nuestros procedimientos son honestos con el cliente para certificar la calidad del producto y la felicidad del cliente que lo compró"

$ws.Range("B15").Value = "Responsabilidad social"
$ws.Range("C15").Value = "Tenemos la obligación con la comunidad de ser una empresa amable, que se preocupa por el cliente y por sus empleados."

$ws.Range("B16").Value = "Calidad"
$ws.Range("C16").Value = "Nuestros productos siempre tienen los estándares de calidad más altos, sin excepciones."

$ws.Range("B17").Value = "Honestidad"
$ws.Range("C17").Value = "Como empresa colombiana estamos en la obigacion de contribuir con un porcentaje de nuestras ganancias al estado. De igual manera prohibimos cualquier accion que nos pueda catalogar como competencia desleal"

# Row heights to accommodate the newly-added wrapped text
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 60

# Scroll / selection state
$ws.Range("A17:C17").Select()
$ws.Application.ActiveWindow.ScrollRow = 2
